# Add a new worksheet "ODI Bowling Extra" (sheetId 5) after the existing
# "ODI Batting Extra" sheet, matching the scraper's new extra-bowling-stats
# export, and populate it with its header row + data rows.

$wb = $excel.ActiveWorkbook

# Place the new sheet immediately after the last existing sheet
# ("ODI Batting Extra"), so the tab order matches the target workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "ODI Bowling Extra"

# Reuse the exact header style (bold, centered, bordered) already used by
# every other sheet's header row (copy format from "ODI Batting Extra"!A1).
$wb.Worksheets.Item("ODI Batting Extra").Range("A1").Copy() | Out-Null
$newSheet.Range("A1:C1").PasteSpecial(-4122) | Out-Null

# Header row
$newSheet.Cells.Item(1,1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1,2).Value = "MAIDEN_OVERS"
$newSheet.Cells.Item(1,3).Value = "PERCENT_WICKETS_OF_ALL"

# Data rows (all values are stored as text, matching the source export;
# the leading "'" forces text storage instead of Excel auto-detecting a
# number/percentage).
$newSheet.Cells.Item(2,1).Value = "'4050"
$newSheet.Cells.Item(3,1).Value = "'4079"
$newSheet.Cells.Item(4,1).Value = "'4081"
$newSheet.Cells.Item(4,2).Value = "'0"
$newSheet.Cells.Item(4,3).Value = "'10.00%"
$newSheet.Cells.Item(5,1).Value = "'4082"
$newSheet.Cells.Item(5,2).Value = "'0"
$newSheet.Cells.Item(5,3).Value = "'10.00%"
$newSheet.Cells.Item(6,1).Value = "'4084"
$newSheet.Cells.Item(6,2).Value = "'1"
$newSheet.Cells.Item(7,1).Value = "'4087"
$newSheet.Cells.Item(7,2).Value = "'0"
$newSheet.Cells.Item(8,1).Value = "'4223"
$newSheet.Cells.Item(8,2).Value = "'0"
$newSheet.Cells.Item(9,1).Value = "'4225"
$newSheet.Cells.Item(10,1).Value = "'4237"
$newSheet.Cells.Item(10,2).Value = "'0"
$newSheet.Cells.Item(11,1).Value = "'4238"
$newSheet.Cells.Item(11,2).Value = "'0"
$newSheet.Cells.Item(12,1).Value = "'4247"
$newSheet.Cells.Item(13,1).Value = "'4297"
$newSheet.Cells.Item(13,2).Value = "'0"
$newSheet.Cells.Item(14,1).Value = "'4300"
$newSheet.Cells.Item(14,2).Value = "'0"
$newSheet.Cells.Item(15,1).Value = "'4308"
$newSheet.Cells.Item(16,1).Value = "'4319"
$newSheet.Cells.Item(16,2).Value = "'0"
$newSheet.Cells.Item(16,3).Value = "'10.00%"
$newSheet.Cells.Item(17,1).Value = "'4324"
$newSheet.Cells.Item(17,2).Value = "'0"
$newSheet.Cells.Item(18,1).Value = "'4334"
$newSheet.Cells.Item(18,2).Value = "'0"
$newSheet.Cells.Item(19,1).Value = "'4337"
$newSheet.Cells.Item(20,1).Value = "'4340"
$newSheet.Cells.Item(20,2).Value = "'0"
$newSheet.Cells.Item(21,1).Value = "'4349"
$newSheet.Cells.Item(21,2).Value = "'1"
